# feat: Weapon 데이터 추가
# Add a new "EWeaponType" enum row (row 3) to the Enum sheet, listing
# the weapon type members in columns A-G.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Enum")

$ws.Range("A3").Value = "EWeaponType"
$ws.Range("B3").Value = "Pistol"
$ws.Range("C3").Value = "Rifle"
$ws.Range("D3").Value = "Shotgun"
$ws.Range("E3").Value = "Sniper"
$ws.Range("F3").Value = "Machinegun"
$ws.Range("G3").Value = "Launcher"

# Match the widened/re-fit columns that Excel produced once the new,
# wider weapon-type values were entered.
$ws.Columns.Item(3).ColumnWidth = 5.375
$ws.Columns.Item(4).ColumnWidth = 8.75
$ws.Columns.Item(5).ColumnWidth = 6.875
$ws.Columns.Item(6).ColumnWidth = 12.5
$ws.Columns.Item(7).ColumnWidth = 9.25

# Leave the selection where the author left it after the edit.
$ws.Range("C5").Select() | Out-Null
